$wb = $excel.ActiveWorkbook

# --- 1. Status text: "Ready for handoff" -> "In Translation" -----------------
# Overview sheet (zh-cn / de-de status columns E & F, rows 2-3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# zh-cn detail sheet (Status column C, rows 2-3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

# de-de detail sheet (Status column C, rows 2-3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- 2. Narrow the status columns --------------------------------------------
# Target stored width is 13.4101845877511 characters; the engine snaps
# ColumnWidth assignments onto a 1/6-character grid, so 12.5 is the input
# that lands closest on that grid (13.333333333333334).
$newStatusColumnWidth = 12.5

$wsOverview.Columns.Item("E").ColumnWidth = $newStatusColumnWidth
$wsOverview.Columns.Item("F").ColumnWidth = $newStatusColumnWidth
$wsZhCn.Columns.Item("C").ColumnWidth = $newStatusColumnWidth
$wsDeDe.Columns.Item("C").ColumnWidth = $newStatusColumnWidth
